# Weekly cryptos-list refresh (GitHub Actions data pull).
# Column D ("Price") cells that are plain decimal numbers get coerced to
# numeric values by Excel unless the cell is pre-formatted as Text - the
# source feed stores every Price/Volume cell as text (note values like
# "27.875.86" that use "." as a thousands separator), so force Text format
# on column D before writing so trailing zeros / separators survive.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.875.86'
$ws.Range('E2').Value = '  -0.21%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.629.48'
$ws.Range('E3').Value = '  -0.99%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.997'
$ws.Range('E4').Value = '  -0.43%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.60'
$ws.Range('E5').Value = '  -0.94%  '
$ws.Range('E6').Value = '  -0.92%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.996'
$ws.Range('E7').Value = '  -0.45%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.25'
$ws.Range('E8').Value = '  -0.71%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.258'
$ws.Range('E9').Value = '  -2.65%  '
$ws.Range('E10').Value = '  -0.26%  '
$ws.Range('E11').Value = '  +0.94%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.858.58'
$ws.Range('E12').Value = '  -1.09%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.627.53'
$ws.Range('E13').Value = '  -1.13%  '
$ws.Range('E14').Value = '  -0.76%  '
$ws.Range('E15').Value = '  -0.40%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.26'
$ws.Range('E16').Value = '  -0.52%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.862.82'
$ws.Range('E17').Value = '  -0.27%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '230.46'
$ws.Range('E18').Value = '  -0.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0722'
$ws.Range('E19').Value = '  -0.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.52'
$ws.Range('E20').Value = '  -1.85%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.996'
$ws.Range('E21').Value = '  -0.50%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.36'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.34'
$ws.Range('E23').Value = '  -3.06%  '
$ws.Range('E24').Value = '  -4.27%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.63'
$ws.Range('E25').Value = '  +1.37%  '
$ws.Range('E26').Value = '  +0.40%  '
$ws.Range('E27').Value = '  -1.07%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.59'
$ws.Range('E28').Value = '  -1.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.997'
$ws.Range('E29').Value = '  -0.48%  '
$ws.Range('E30').Value = '  -1.11%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0481'
$ws.Range('E31').Value = '  -0.93%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.41'
$ws.Range('E32').Value = '  +2.19%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.406.27'
$ws.Range('E33').Value = '  -2.54%  '
$ws.Range('E34').Value = '  +0.10%  '
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('E36').Value = '  +8.60%  '
$ws.Range('E37').Value = '  +0.90%  '
$ws.Range('E38').Value = '  +0.44%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.558'
$ws.Range('E39').Value = '  -0.08%  '
$ws.Range('E40').Value = '  -2.63%  '
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.996'
$ws.Range('E42').Value = '  -0.48%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '66.58'
$ws.Range('E43').Value = '  -3.71%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.51'
$ws.Range('E44').Value = '  +1.54%  '
$ws.Range('E45').Value = '  -0.23%  '
$ws.Range('E46').Value = '  -1.11%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.769.28'
$ws.Range('E47').Value = '  -1.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '87.81'
$ws.Range('E48').Value = '  -1.44%  '

# Row 49: a new listing (BabyDogeCoin) now leads the tail of the table, pushing
# Algorand -> row 50 and Cronos -> row 51; EnergySwap (old row 51) drops off
# the bottom of the (still 51-row) list. Column A's rank index is untouched.
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0₆0105'
$ws.Range('E49').Value = '  +0.01%  '

$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0997'
$ws.Range('E50').Value = '  -1.15%  '

$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0506'
$ws.Range('E51').Value = '  -0.39%  '
